$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.530.51'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '1.664.16'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.45'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4808'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2632'
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06168'
$ws.Range("E9").Value = '  +2.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07110'
$ws.Range("E10").Value = '  -1.05%  '
$ws.Range("D11").Value = '1.662.78'
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.79'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5887'
$ws.Range("E13").Value = '  -5.34%  '
$ws.Range("E14").Value = '  -4.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.08'
$ws.Range("E15").Value = '  +2.56%  '
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9998'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").Value = '25.519.98'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006744'
$ws.Range("E19").Value = '  +1.81%  '
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '1.872.55'
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.439'
$ws.Range("E22").Value = '  -2.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.710'
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.289'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.82'
$ws.Range("E25").Value = '  +2.29%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.382'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '105.19'
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.709'
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.956'
$ws.Range("E30").Value = '  +4.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.657'
$ws.Range("E31").Value = '  +1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07723'
$ws.Range("E32").Value = '  -2.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9991'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04247'
$ws.Range("E34").Value = '  -7.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.601'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6120'
$ws.Range("E36").Value = '  +6.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9498'
$ws.Range("E37").Value = '  +0.79%  '
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8620'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9994'
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.856'
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01467'
$ws.Range("E42").Value = '  -6.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.93'
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3766'
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.844'
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1123'
$ws.Range("E46").Value = '  -1.88%  '
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05259'
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.75'
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("E51").Value = '  -0.03%  '
